# Loan RBI, Variable Instalments
# Insert a new (blank) column into the "Repayment schedule" sheet so that a
# "Late" / "Outstanding" pair of columns is split apart by a new, currently
# empty "Variable Instalments" style column, and make the "Repayment
# schedule" sheet the active tab (it was "NewLoanInput" before).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N..P to O..Q)
$ws.Columns("N").Insert()
$ws.Columns("N").ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab and update the selected
# cell to match the authored workbook state.
$ws.Select()
$ws.Range("K18").Select() | Out-Null
